# WASHINGTON_2017.xlsx data-cleaning fix
# - Rename header columns to short machine-friendly codes
# - Title-case Spanish connector words (de/del/la/las/el/los/y) in the
#   "Estado de Origen" / "Municipio Origen" text columns
# - Remove the trailing free-text metadata/footnote rows
# - Tidy a floating point rounding artifact in D566

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Convert-SpanishTitleCase($text) {
    $result = $text
    $result = $result -creplace '(?<=^| )de(?=$| )', 'De'
    $result = $result -creplace '(?<=^| )del(?=$| )', 'Del'
    $result = $result -creplace '(?<=^| )la(?=$| )', 'La'
    $result = $result -creplace '(?<=^| )las(?=$| )', 'Las'
    $result = $result -creplace '(?<=^| )el(?=$| )', 'El'
    $result = $result -creplace '(?<=^| )los(?=$| )', 'Los'
    $result = $result -creplace '(?<=^| )y(?=$| )', 'Y'
    return $result
}

# 1) Apply the Spanish title-case cleanup to columns A (state) and B (municipality)
#    for every data row (skip the header row 1).
for ($r = 2; $r -le 1355; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string]) {
            $newVal = Convert-SpanishTitleCase $val
            if (-not ($newVal.Equals($val))) {
                $cell.Value = $newVal
            }
        }
    }
}

# 2) Rename the header row to short machine-friendly field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Minor float-precision tidy-up on D566
$ws.Cells.Item(566, 4).Value2 = 0.009895337773549

# 4) Drop the trailing free-text footer/metadata rows (sample size, source,
#    author, etc.) that lived below the data table.
$ws.Range("A1357:D1361").EntireRow.Delete() | Out-Null

Write-Host "Done. New dimension:" $ws.UsedRange.Address()
